$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the QuestDungeonRate for the row with Id 18000201 (row 6):
# add a new drop entry "safebox;1" to the existing "bookancient;1" value.
$ws.Range("M6").Value = "bookancient;1|safebox;1"

# Reflect the new selection left after the edit.
$ws.Range("M6").Select()
